# Commit: "Fruta / hortaliza, semanal" — weekly refresh of the Hortaliza /
# Vega Monumental Concepción - Zanahoria dataset: two new daily records are
# inserted at the top of the historical block (rows 189-190), pushing the
# rest of the series down by two rows (old 189..269 -> new 191..271).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows right before the current row 189 (done one at a
# time, mirroring Excel's native Rows.Insert behaviour — each insert shifts
# everything at/below row 189 down by one row).
$ws.Rows.Item(189).EntireRow.Insert()
$ws.Rows.Item(189).EntireRow.Insert()

# New record #1 -> row 189
$ws.Range("A189").Value = 11
$ws.Range("B189").Value = "Vega Monumental Concepción"
$ws.Range("C189").Value = "Bíobío"
$ws.Range("D189").Value = 44813
$ws.Range("E189").Value = 8
$ws.Range("F189").Value = 100114013
$ws.Range("G189").Value = "Zanahoria"
$ws.Range("H189").Value = "Sin especificar"
$ws.Range("I189").Value = "Primera"
$ws.Range("J189").Value = 600
$ws.Range("K189").Value = 8000
$ws.Range("L189").Value = 9000
$ws.Range("M189").Value = 8500
$ws.Range("N189").Value = "`$/saco 20 kilos"
$ws.Range("O189").Value = "Región de Ñuble"
$ws.Range("P189").Value = 425
$ws.Range("Q189").Value = 20
$ws.Range("R189").Value = "Hortaliza"

# New record #2 -> row 190
$ws.Range("A190").Value = 11
$ws.Range("B190").Value = "Vega Monumental Concepción"
$ws.Range("C190").Value = "Bíobío"
$ws.Range("D190").Value = 44813
$ws.Range("E190").Value = 8
$ws.Range("F190").Value = 100114013
$ws.Range("G190").Value = "Zanahoria"
$ws.Range("H190").Value = "Sin especificar"
$ws.Range("I190").Value = "Segunda"
$ws.Range("J190").Value = 300
$ws.Range("K190").Value = 7000
$ws.Range("L190").Value = 7000
$ws.Range("M190").Value = 7000
$ws.Range("N190").Value = "`$/saco 20 kilos"
$ws.Range("O190").Value = "Región de Ñuble"
$ws.Range("P190").Value = 350
$ws.Range("Q190").Value = 20
$ws.Range("R190").Value = "Hortaliza"
